$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "mngr250218"
$ws.Range("B2").Value = "vezYgad"
